$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to coerce numeric-looking strings ("32", "49%", "705", ...)
# into genuine text values without Excel's automatic number inference, and
# without leaving behind any unused NumberFormat/style entries.
$scratch = "Z100"

function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $ws.Range($scratch).Formula = '="' + $escaped + '"'
    $ws.Range($range).Formula = "=T(" + $scratch + ")"
    $ws.Range($range).Copy()
    $ws.Range($range).PasteSpecial(-4163)
    $ws.Range($scratch).Clear()
}

# ------------------------------------------------------------------
# Header row (row 1): extend B1:D1's existing bordered/bold style
# across the new header cells E1:S1, then fill in the header text.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("E1:S1").PasteSpecial(-4122)

$headerCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
$headers = @(
    "name","nickname","record","status","hometown","trains_at","age","height",
    "weight","octagon_debut","reach","leg_reach","significant_strikes_landed",
    "significant_strikes_attempted","signicant_strike_accuracy","takedowns_landed",
    "takedowns_attempted","takedowns_accuracy"
)
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $col = $headerCols[$i]
    if ($col -eq "B" -or $col -eq "C" -or $col -eq "D") {
        continue # already populated/styled in the original workbook
    }
    $ws.Range($col + "1").Value = $headers[$i]
}

# ------------------------------------------------------------------
# Row 3: give A3 the same bordered "index" style as A2, then set the
# numeric record-index values for both rows.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1

# ------------------------------------------------------------------
# Data rows 2 (Khabib Nurmagomedov) and 3 (Jon Jones), column by
# column, writing row 2's value then row 3's value - this mirrors the
# order the fighter stats were originally scraped/appended in.
# ------------------------------------------------------------------
$row2 = @{
    "B" = "Khabib Nurmagomedov"
    "C" = "The Eagle"
    "D" = "29-0-0 (W-L-D)"
    "E" = "Retired"
    "F" = "Dagestan Republic, Russia"
    "G" = "AKA (American Kickboxing Academy) San Jose"
    "H" = "32"
    "I" = "70.00"
    "J" = "155.00"
    "K" = "Jan. 21, 2012"
    "L" = "70.00"
    "M" = "40.00"
    "N" = "705"
    "O" = "1444"
    "P" = "49%"
    "Q" = "49"
    "R" = "127"
    "S" = "48%"
}
$row3 = @{
    "B" = "Jon Jones"
    "C" = "Bones"
    "D" = "26-1-0 (W-L-D)"
    "E" = "Active"
    "F" = "Rochester, United States"
    "H" = "33"
    "I" = "76.00"
    "J" = "205.00"
    "K" = "Aug. 09, 2008"
    "L" = "84.50"
    "M" = "45.00"
    "N" = "1463"
    "O" = "2526"
    "P" = "58%"
    "Q" = "36"
    "R" = "95"
    "S" = "44%"
}

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
foreach ($col in $cols) {
    if ($row2.ContainsKey($col)) {
        Set-TextValue ($col + "2") $row2[$col]
    }
    if ($row3.ContainsKey($col)) {
        Set-TextValue ($col + "3") $row3[$col]
    }
}
